$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D29").Value = "프로메디우스"

$ws.Range("D37").Value = "[Paper Review] ConViT : Improving Vision Transformers with Soft Convolutional Inductive Biases"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1842&mod=document&pageid=1"

$ws.Range("D46").Value = "[씨젠] 2021년 11월, 생물정보학(Bioinformatics 채용), Bioinformatics Data Engineer, Bioinformatics SW Engineer"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/421"
